$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# NOTE: in this engine, Range.Hyperlinks.Delete() clears *every* hyperlink on
# the worksheet (not just the target range), so do it once up front and
# rebuild the surviving links (F4, F8, F9) afterward.
$ws.Range("F2").Hyperlinks.Delete()

# ---- Row 2: Utilite Standard -> BeagleBone Black ----
$ws.Range("B2").Value = "BeagleBone Black"
$ws.Range("C2").Value = 32
$ws.Range("D2").Value = 55
$ws.Range("F2").Value = "http://www.mouser.com/ProductDetail/BeagleBoard-by-CircuitCo/BB-BBLK-000/?qs=%2fha2pyFadugh6wNMONnDuAbTwbrIHVw4R%2f%252bth5Q2M%2fX2Gs60muroNw%3d%3d"

# ---- Row 3: HP 2920-24G -> HP 2920-48G (plain text now, no live hyperlink) ----
$ws.Range("B3").Value = "HP 2920-48G Network Switch"
$ws.Range("F3").Value = "http://www.amazon.com/HP-J9728A-2920-48G-Switch/dp/B00BJ42JQY"

# ---- Row 4: Kingston MicroSD cards qty 16 -> 32 (url/name unchanged) ----
$ws.Range("C4").Value = 32

# ---- Row 5: clear part details (istarUSA WN228 rack chassis removed) ----
$ws.Range("B5:D5").Clear()
$ws.Range("F5").ClearContents()

# ---- Row 6: clear part details (istarUSA WA-PS010 removed) ----
$ws.Range("B6:D6").Clear()
$ws.Range("F6").ClearContents()

# ---- Row 7: istarUSA WA-SF80B flat vented tray -> 10/100 Ethernet Adapter (w/ linux) ----
$ws.Range("B7").Value = "10/100 Ethernet Adapter (w/ linux)"
$ws.Range("C7").Value = 32
$ws.Range("D7").Value = 11.99
$ws.Range("F7").Value = "http://www.newegg.com/Product/Product.aspx?Item=N82E16812315001&cm_re=usb_ethernet-_-12-315-001-_-Product"

# ---- Row 8: Ethernet Cable (5-set) qty 10 -> 15 (url/name unchanged) ----
$ws.Range("C8").Value = 15

# ---- Row 9: TRENDnet switch - unchanged ----

# ---- Re-create the hyperlinks that survive in the final sheet ----
$ws.Hyperlinks.Add($ws.Range("F4"), "http://www.amazon.com/Kingston-Digital-microSDHC-SDC4-16GBET/dp/B00DYQYLQQ/ref=sr_1_1?ie=UTF8&qid=1422292553&sr=8-1&keywords=micro+sd+16gb") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "http://www.amazon.com/Cable-Matters-5-Color-Snagless-Ethernet/dp/B00E5I7VJG/ref=sr_1_1?s=pc&ie=UTF8&qid=1422383476&sr=1-1&keywords=ethernet+cables") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "http://www.amazon.com/TRENDnet-24-Port-Unmanaged-GREENnet-TEG-S24Dg/dp/B0044GH27U/ref=sr_1_4?s=electronics&ie") | Out-Null

# Hyperlinks.Add() stamps a duplicate "Hyperlink" style record onto its
# target cell instead of reusing the sheet's existing one (style index 7,
# already applied to every cell in column F) - copy/paste the formatting
# from an untouched F-column cell so F4/F8/F9 land back on that same index.
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- New rows 11 & 12 ----
$ws.Range("B11").Value = "min cables required"
$ws.Range("C11").Formula = "=C2*2"

$ws.Range("B12").Value = "cables ordered"
$ws.Range("C12").Formula = "=C8*5"

$ws.Range("D14").Select() | Out-Null
